# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) timestamps for the zh-cn and de-de
# handback-status sheets to reflect the newly generated report times.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2:D3").Value = "2016-02-18 04:17:17"
$wsZhCn.Range("G2:G3").Value = "2016-02-18 04:18:05"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2:D3").Value = "2016-02-18 04:17:30"
$wsDeDe.Range("G2:G3").Value = "2016-02-18 04:18:29"
